$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 13 (pushes old "Programa resumido:" block and everything after down by 2)
$ws.Range("A13:A14").EntireRow.Insert()

# Fix / fill in B & C cell values that changed position or content
$ws.Range("B10").Value = 'Propiciar conhecimentos teóricos e práticos de materiais e ferramentas abrasivas, incluindo as matérias primas, os processos de fabricação, qualificação, teste e aplicação de ferramentas abrasivas; e os mecanismos de abrasão atuantes durante os processos de usinagem. Aprimorar a formação do Engenheiro de Materiais numa área estratégica de Ciência e Engenharia de Materiais presente em diversos segmentos da indústria de transformação do país'
$ws.Range("C10").Value = 'Propiciar conhecimentos teóricos e práticos de materiais e ferramentas abrasivas, incluindo as matérias primas, os processos de fabricação, qualificação, teste e aplicação de ferramentas abrasivas; e os mecanismos de abrasão atuantes durante os processos de usinagem. Aprimorar a formação do Engenheiro de Materiais numa área estratégica de Ciência e Engenharia de Materiais presente em diversos segmentos da indústria de transformação do país'

$ws.Range("B13").Value = '519033 - Carlos Yujiro Shigue'
$ws.Range("C13").Value = '519033 - Carlos Yujiro Shigue'

$ws.Range("B14").Value = '5817692 - Katia Cristiane Gandolpho Candioto'
$ws.Range("C14").Value = '5817692 - Katia Cristiane Gandolpho Candioto'

$ws.Range("B15").Value = 'Introdução aos materiais e ferramentas abrasivas. Matérias primas utilizadas na fabricação de ferramentas abrasivas. Processos de fabricação de ferramentas abrasivas. Caracterização, teste e inspeção de ferramentas abrasivas. Mecânica da usinagem com ferramentas abrasivas. Avaliação de desempenho. Prática laboratorial de fabricação, caracterização e teste de compósitos abrasivos.'
$ws.Range("C15").Value = 'Introdução aos materiais e ferramentas abrasivas. Matérias primas utilizadas na fabricação de ferramentas abrasivas. Processos de fabricação de ferramentas abrasivas. Caracterização, teste e inspeção de ferramentas abrasivas. Mecânica da usinagem com ferramentas abrasivas. Avaliação de desempenho. Prática laboratorial de fabricação, caracterização e teste de compósitos abrasivos.'

$ws.Range("B17").Value = 'Conteúdo teórico: 1. Introdução aos materiais e ferramentas abrasivas: histórico, materiais abrasivos naturais e sintéticos e características principais dos materiais abrasivos. Dados econômicos das ferramentas abrasivas.2. Matérias primas utilizadas na fabricação de ferramentas abrasivas: cerâmicas, borrachas, metais e polímeros. Processos de obtenção das matérias primas.3. Processos de fabricação de ferramentas abrasivas: discos, rebolos, pontas montadas e lixas.4. Caracterização, teste e inspeção de ferramentas abrasivas: ensaios destrutivos e não destrutivos. Normas e códigos de segurança. 5. Mecânica da usinagem com ferramentas abrasivas. Operações com abrasivos: corte, retificação, desbaste, acabamento, lapidação e afiação.6. Avaliação de desempenho: aspectos térmicos, refrigeração, lubrificação, rugosidade superficial, interação metal-ferramenta e defeitos em ferramentas abrasivas.Conteúdo prático: 1. Prática laboratorial de fabricação, caracterização e teste de compósitos abrasivos.2. Visita a fabricantes de abrasivos. 3. Visita a usuários de ferramentas abrasivas.'
$ws.Range("C17").Value = 'Conteúdo teórico: 1. Introdução aos materiais e ferramentas abrasivas: histórico, materiais abrasivos naturais e sintéticos e características principais dos materiais abrasivos. Dados econômicos das ferramentas abrasivas.2. Matérias primas utilizadas na fabricação de ferramentas abrasivas: cerâmicas, borrachas, metais e polímeros. Processos de obtenção das matérias primas.3. Processos de fabricação de ferramentas abrasivas: discos, rebolos, pontas montadas e lixas.4. Caracterização, teste e inspeção de ferramentas abrasivas: ensaios destrutivos e não destrutivos. Normas e códigos de segurança. 5. Mecânica da usinagem com ferramentas abrasivas. Operações com abrasivos: corte, retificação, desbaste, acabamento, lapidação e afiação.6. Avaliação de desempenho: aspectos térmicos, refrigeração, lubrificação, rugosidade superficial, interação metal-ferramenta e defeitos em ferramentas abrasivas.Conteúdo prático: 1. Prática laboratorial de fabricação, caracterização e teste de compósitos abrasivos.2. Visita a fabricantes de abrasivos. 3. Visita a usuários de ferramentas abrasivas.'

$ws.Range("B20").Value = 'A avaliação será constituída por aulas expositivas, aulas de exercícios e práticas laboratoriais. Serão aplicadas pelo menos duas avaliações.'
$ws.Range("C20").Value = 'A avaliação será constituída por aulas expositivas, aulas de exercícios e práticas laboratoriais. Serão aplicadas pelo menos duas avaliações.'

$ws.Range("B21").Value = 'A nota final será a média das avaliações escritas e práticas'
$ws.Range("C21").Value = 'A nota final será a média das avaliações escritas e práticas'

$ws.Range("B22").Value = 'A recuperação será uma prova escrita (RE) que comporá com a nota final (NF) a média final (MF), sendo MF = (NF + RE)/2.'
$ws.Range("C22").Value = 'A recuperação será uma prova escrita (RE) que comporá com a nota final (NF) a média final (MF), sendo MF = (NF + RE)/2.'

$ws.Range("B23").Value = '1. NUSSBAUM, G. C. Rebolos e abrasivos. Tecnologia básica. São Paulo: Ícone Editora, 1988. 2. KLOCKE, F. Manufacturing processes 2. Grinding, honing, lapping. Berlim: Springer Verlag, 2009.3. MALKIN, S.; GUO, C. Grinding technology: theory and application of machining with abrasives. New York: Industrial Press Inc., 2008.4. JACKSON, M. J.; DAVIM, J. P. Machining with abrasives. New York: Springer Science, 2011.5. FERRARESI, D. Usinagem dos metais. São Paulo: Editora Edgard Blucher, 1970.6. STEMMER, C. E. Ferramentas de corte II: brocas, alargadores, ferramentas de rocar, fresas, brochas, rebolos e abrasivos. Florianópolis: Editora da UFSC, 1995.7. KINGERY, W. D. Ceramic fabrication process. New York: John Wiley, 1958.8. GARDZIELLA, A.; PILATO, L.A.; KNOP, A. Phenolic resins: chemistry, applications, standardization, safety and ecology. Berlim: Springer Verlag, 2000.9. MARINESCU, Ioan D. Tribology of abrasive machining processes. 2ª Ed. New York: Willian Andrew, 2004.'
$ws.Range("C23").Value = '1. NUSSBAUM, G. C. Rebolos e abrasivos. Tecnologia básica. São Paulo: Ícone Editora, 1988. 2. KLOCKE, F. Manufacturing processes 2. Grinding, honing, lapping. Berlim: Springer Verlag, 2009.3. MALKIN, S.; GUO, C. Grinding technology: theory and application of machining with abrasives. New York: Industrial Press Inc., 2008.4. JACKSON, M. J.; DAVIM, J. P. Machining with abrasives. New York: Springer Science, 2011.5. FERRARESI, D. Usinagem dos metais. São Paulo: Editora Edgard Blucher, 1970.6. STEMMER, C. E. Ferramentas de corte II: brocas, alargadores, ferramentas de rocar, fresas, brochas, rebolos e abrasivos. Florianópolis: Editora da UFSC, 1995.7. KINGERY, W. D. Ceramic fabrication process. New York: John Wiley, 1958.8. GARDZIELLA, A.; PILATO, L.A.; KNOP, A. Phenolic resins: chemistry, applications, standardization, safety and ecology. Berlim: Springer Verlag, 2000.9. MARINESCU, Ioan D. Tribology of abrasive machining processes. 2ª Ed. New York: Willian Andrew, 2004.'
